$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 30 cell update(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3916.1428
$ws.Range("I62").Value = 3916.1428
$ws.Range("K62").Value = 3916.1428
$ws.Range("M62").Value = -3292.1428
$ws.Range("H65").Value = 3916.1428
$ws.Range("I65").Value = 3916.1428
$ws.Range("K65").Value = 19580.714
$ws.Range("M65").Value = -16460.714
$ws.Range("H74").Value = 9473.684999999999
$ws.Range("I74").Value = 9722.223
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 9722.223
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -8786.223
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 9473.684999999999
$ws.Range("I77").Value = 9722.223
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 48611.115
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -43931.115
$ws.Range("N77").Value = -34360
$ws.Range("H137").Value = 4079.9524
$ws.Range("I137").Value = 1719.1621
$ws.Range("K137").Value = 5157.4863
$ws.Range("M137").Value = -2607.4863
$ws.Range("H138").Value = 3269.8928
$ws.Range("I138").Value = 2512.238
$ws.Range("K138").Value = 7536.714
$ws.Range("M138").Value = -2396.714

# --- Sheet ARM: 38 cell update(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 190334.58
$ws.Range("I32").Value = 195655.14
$ws.Range("K32").Value = 195655.14
$ws.Range("M32").Value = -195368.14
$ws.Range("H61").Value = 3129.3333
$ws.Range("I61").Value = 3181.4285
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 3181.4285
$ws.Range("L61").Value = 2400
$ws.Range("M61").Value = -2969.4285
$ws.Range("N61").Value = -2824
$ws.Range("H74").Value = 6496.735
$ws.Range("I74").Value = 3923.7144
$ws.Range("K74").Value = 3923.7144
$ws.Range("M74").Value = -3049.7144
$ws.Range("H77").Value = 6496.735
$ws.Range("I77").Value = 3923.7144
$ws.Range("K77").Value = 19618.572
$ws.Range("M77").Value = -15250.572
$ws.Range("H110").Value = 2029.2354
$ws.Range("I110").Value = 2083
$ws.Range("K110").Value = 2083
$ws.Range("M110").Value = -38
$ws.Range("H132").Value = 4260.0205
$ws.Range("I132").Value = 2520.861
$ws.Range("K132").Value = 7562.583
$ws.Range("M132").Value = -5032.583
$ws.Range("H136").Value = 3129.3333
$ws.Range("I136").Value = 3181.4285
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 9544.2855
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -6994.2855
$ws.Range("N136").Value = -12300
$ws.Range("H139").Value = 85446.875
$ws.Range("J139").Value = 85446.875
$ws.Range("L139").Value = 85446.875
$ws.Range("N139").Value = -95726.875

# --- Sheet BSM: 22 cell update(s) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21680.5
$ws.Range("J82").Value = 55000
$ws.Range("L82").Value = 55000
$ws.Range("N82").Value = -55766
$ws.Range("H85").Value = 21680.5
$ws.Range("J85").Value = 55000
$ws.Range("L85").Value = 55000
$ws.Range("N85").Value = -57652
$ws.Range("H94").Value = 5429
$ws.Range("I94").Value = 4781.467
$ws.Range("J94").Value = 8666.666999999999
$ws.Range("K94").Value = 4781.467
$ws.Range("L94").Value = 8666.666999999999
$ws.Range("M94").Value = -4330.467
$ws.Range("N94").Value = -9568.666999999999
$ws.Range("H134").Value = 5769.2915
$ws.Range("I134").Value = 6123.5
$ws.Range("J134").Value = 3998.25
$ws.Range("K134").Value = 18370.5
$ws.Range("L134").Value = 11994.75
$ws.Range("M134").Value = -15835.5
$ws.Range("N134").Value = -17064.75

# --- Sheet CRP: 23 cell update(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3309.7317
$ws.Range("I58").Value = 2713.8076
$ws.Range("K58").Value = 2713.8076
$ws.Range("M58").Value = -2510.8076
$ws.Range("H86").Value = 79360.38
$ws.Range("I86").Value = 131684
$ws.Range("K86").Value = 131684
$ws.Range("M86").Value = -130561
$ws.Range("H89").Value = 79360.38
$ws.Range("I89").Value = 131684
$ws.Range("K89").Value = 658420
$ws.Range("M89").Value = -652804
$ws.Range("H122").Value = 10583.533
$ws.Range("I122").Value = 2168.28
$ws.Range("J122").Value = 52659.8
$ws.Range("K122").Value = 6504.84
$ws.Range("L122").Value = 157979.4
$ws.Range("M122").Value = -4054.84
$ws.Range("N122").Value = -162879.4
$ws.Range("H136").Value = 3309.7317
$ws.Range("I136").Value = 2713.8076
$ws.Range("K136").Value = 8141.4228
$ws.Range("M136").Value = -5591.4228

# --- Sheet CUL: 25 cell update(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 759.5714
$ws.Range("I92").Value = 368
$ws.Range("J92").Value = 1281.6666
$ws.Range("K92").Value = 1104
$ws.Range("L92").Value = 3844.9998
$ws.Range("M92").Value = 144
$ws.Range("N92").Value = -6340.9998
$ws.Range("H97").Value = 775.75
$ws.Range("I97").Value = 562
$ws.Range("K97").Value = 1686
$ws.Range("M97").Value = -1190
$ws.Range("H107").Value = 4134.143
$ws.Range("I107").Value = 2178.3333
$ws.Range("J107").Value = 4667.5454
$ws.Range("K107").Value = 6534.999899999999
$ws.Range("L107").Value = 14002.6362
$ws.Range("M107").Value = -4614.999899999999
$ws.Range("N107").Value = -17842.6362
$ws.Range("H120").Value = 13333.333
$ws.Range("I120").Value = 13333.333
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 39999.999
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -35161.999
$ws.Range("N120").ClearContents()

# --- Sheet GSM: 27 cell update(s) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 7750
$ws.Range("I3").Value = 9750
$ws.Range("K3").Value = 9750
$ws.Range("M3").Value = -9634
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("H107").Value = 1624
$ws.Range("I107").Value = 1624
$ws.Range("K107").Value = 1624
$ws.Range("M107").Value = 296
$ws.Range("H122").Value = 2619.0476
$ws.Range("I122").Value = 1711.1111
$ws.Range("J122").Value = 3300
$ws.Range("K122").Value = 5133.3333
$ws.Range("L122").Value = 9900
$ws.Range("M122").Value = -2683.3333
$ws.Range("N122").Value = -14800
$ws.Range("H126").Value = 2747.76
$ws.Range("I126").Value = 2545.5293
$ws.Range("K126").Value = 7636.5879
$ws.Range("M126").Value = -5166.5879
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

# --- Sheet LTW: 27 cell update(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2479.3809
$ws.Range("J22").Value = 3346.625
$ws.Range("L22").Value = 3346.625
$ws.Range("N22").Value = -3936.625
$ws.Range("H27").Value = 2479.3809
$ws.Range("J27").Value = 3346.625
$ws.Range("L27").Value = 3346.625
$ws.Range("N27").Value = -3560.625
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40912
$ws.Range("H46").Value = 4242.6313
$ws.Range("I46").Value = 1154.75
$ws.Range("K46").Value = 1154.75
$ws.Range("M46").Value = -966.75
$ws.Range("H100").Value = 2931.5
$ws.Range("I100").Value = 3008.5
$ws.Range("J100").Value = 2777.5
$ws.Range("K100").Value = 3008.5
$ws.Range("L100").Value = 2777.5
$ws.Range("M100").Value = -2467.5
$ws.Range("N100").Value = -3859.5
$ws.Range("H122").Value = 3911.3333
$ws.Range("I122").Value = 3911.3333
$ws.Range("K122").Value = 11733.9999
$ws.Range("M122").Value = -9283.999899999999

# --- Sheet WVR: 12 cell update(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7846.75
$ws.Range("J81").Value = 17170.143
$ws.Range("L81").Value = 34340.286
$ws.Range("N81").Value = -36462.286
$ws.Range("H84").Value = 7846.75
$ws.Range("J84").Value = 17170.143
$ws.Range("L84").Value = 171701.43
$ws.Range("N84").Value = -182309.43
$ws.Range("H136").Value = 6946
$ws.Range("I136").Value = 7403.385
$ws.Range("K136").Value = 22210.155
$ws.Range("M136").Value = -19660.155
